# Micropipette Presentation workbook edit
# - Rename Sheet1 -> "Drift Test", Sheet2 -> "200uL"
# - Add two new sheets "20uL" and "50uL" with precision/accuracy data
# - Make "200uL" the active tab
# - Update selections on "Drift Test" and the new data sheets

$wb = $excel.ActiveWorkbook

$wsDrift = $wb.Worksheets.Item(1)
$ws200   = $wb.Worksheets.Item(2)

$wsDrift.Name = "Drift Test"
$ws200.Name   = "200uL"

# Add the two new sheets right after "200uL", in order: 20uL, then 50uL
$ws20 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws200)
$ws20.Name = "20uL"
$ws50 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws20)
$ws50.Name = "50uL"

# ===================== 200uL =====================
$ws200.Range("B1").Value = "Test"
$ws200.Range("C1").Formula = "=0.2*0.99997"

$a200 = @(0.2028, 0.2044, 0.2014, 0.2004, 0.201, 0.2052, 0.2082, 0.2065, 0.2075, 0.2072)
$b200 = @(0.1929, 0.1928, 0.1922, 0.2072, 0.1994, 0.1924, 0.2001, 0.1997, 0.1995, 0.1962)
for ($i = 0; $i -lt $a200.Length; $i++) {
    $row = $i + 2
    $ws200.Cells.Item($row, 1).Value = $a200[$i]
    $ws200.Cells.Item($row, 2).Value = $b200[$i]
}

$ws200.Range("A12").Formula = "=AVERAGE(A2:A11)"
$ws200.Range("B12").Formula = "=AVERAGE(B2:B11)"
$ws200.Range("C12").Value = "Average"

$ws200.Range("A13").Formula = "=STDEV(A2:A11)"
$ws200.Range("B13").Formula = "=STDEV(B2:B11)"
$ws200.Range("C13").Value = "Standard Deviation"

$ws200.Range("A14").Formula = "=((A12-(C1))/(C1))*100"
$ws200.Range("B14").Formula = "=(B12-C1)/(C1)*100"
$ws200.Range("C14").Value = "Percent Error"

$ws200.Range("A15").Formula = "=A13*100/A12"
$ws200.Range("B15").Formula = "=B13*100/B12"
$ws200.Range("C15").Value = "Percent Standard Deviation"

# A1 is set last so "Test" (B1) registers in the shared string table
# before "Commercial "
$ws200.Range("A1").Value = "Commercial "

$ws200.Columns.Item(1).ColumnWidth = 12.25

# ===================== 20uL =====================
$ws20.Range("B1").Value = "Test"
$ws20.Range("C1").Formula = "=0.02*0.99997"

$a20 = @(0.0198, 0.0196, 0.02, 0.02, 0.0199, 0.0196, 0.0198, 0.0196, 0.02, 0.02)
$b20 = @(0.0192, 0.019, 0.0184, 0.0187, 0.019, 0.0195, 0.0196, 0.0195, 0.0196, 0.01956)
for ($i = 0; $i -lt $a20.Length; $i++) {
    $row = $i + 2
    $ws20.Cells.Item($row, 1).Value = $a20[$i]
    $ws20.Cells.Item($row, 2).Value = $b20[$i]
}

$ws20.Range("A12").Formula = "=AVERAGE(A2:A11)"
$ws20.Range("B12").Formula = "=AVERAGE(B2:B11)"
$ws20.Range("C12").Value = "Average"

$ws20.Range("A13").Formula = "=STDEV(A2:A11)"
$ws20.Range("B13").Formula = "=STDEV(B2:B11)"
$ws20.Range("C13").Value = "Standard Deviation"

$ws20.Range("A14").Formula = "=((A12-(C1))/(C1))*100"
$ws20.Range("B14").Formula = "=(B12-C1)/(C1)*100"
$ws20.Range("C14").Value = "Percent Error"

$ws20.Range("A15").Formula = "=A13*100/A12"
$ws20.Range("B15").Formula = "=B13*100/B12"
$ws20.Range("C15").Value = "Percent Standard Deviation"

$ws20.Range("A1").Value = "Commercial "

$ws20.Columns.Item(1).ColumnWidth = 14.25

# ===================== 50uL =====================
$ws50.Range("B1").Value = "Test"
$ws50.Range("C1").Formula = "=0.05*0.99997"

$a50 = @(0.0496, 0.0498, 0.0498, 0.05, 0.0501, 0.05, 0.0499, 0.0499, 0.0499, 0.05)
$b50 = @(0.0571, 0.0546, 0.0535, 0.054, 0.0547, 0.0516, 0.0514, 0.0517, 0.0516, 0.0513)
for ($i = 0; $i -lt $a50.Length; $i++) {
    $row = $i + 2
    $ws50.Cells.Item($row, 1).Value = $a50[$i]
    $ws50.Cells.Item($row, 2).Value = $b50[$i]
}

$ws50.Range("A12").Formula = "=AVERAGE(A2:A11)"
$ws50.Range("B12").Formula = "=AVERAGE(B2:B11)"
$ws50.Range("C12").Value = "Average"

$ws50.Range("A13").Formula = "=STDEV(A2:A11)"
$ws50.Range("B13").Formula = "=STDEV(B2:B11)"
$ws50.Range("C13").Value = "Standard Deviation"

$ws50.Range("A14").Formula = "=((A12-(C1))/(C1))*100"
$ws50.Range("B14").Formula = "=(B12-C1)/(C1)*100"
$ws50.Range("C14").Value = "Percent Error"

$ws50.Range("A15").Formula = "=A13*100/A12"
$ws50.Range("B15").Formula = "=B13*100/B12"
$ws50.Range("C15").Value = "Percent Standard Deviation"

$ws50.Range("A1").Value = "Commercial "

# ===================== View / selection adjustments =====================
$wsDrift.Activate()
$wsDrift.Range("A69").Select()

$ws20.Activate()
$ws20.Range("A13").Select()

$ws50.Activate()
$ws50.Range("A13").Select()

$ws200.Activate()
$ws200.Range("A13").Select()

Write-Host "Edit complete"
